# Updates the LR-pairs data table (A2:T9) per Dr Hou's revised NATMI analysis.
# Each sending cluster (ECs, FAPs, M2, sCs) now reports results for BOTH target
# clusters (M2 and sCs) instead of just M2, so the 4-row table becomes 8 rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 8,20
$data[0,0] = "ECs"
$data[0,1] = "Col3a1"
$data[0,2] = "Mag"
$data[0,3] = "M2"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 21.62966433333333
$data[0,7] = 64.888993
$data[0,8] = 0.004276908378962984
$data[0,9] = 0.004276908378962984
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 0.7885686666666668
$data[0,13] = 2.365706
$data[0,14] = 0.4566863346753138
$data[0,15] = 0.4566863346753137
$data[0,16] = 17.05647556378423
$data[0,17] = 153.508280074058
$data[0,18] = 0.001953205611330743
$data[0,19] = 0.001953205611330742
$data[1,0] = "ECs"
$data[1,1] = "Col3a1"
$data[1,2] = "Mag"
$data[1,3] = "sCs"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 21.62966433333333
$data[1,7] = 64.888993
$data[1,8] = 0.004276908378962984
$data[1,9] = 0.004276908378962984
$data[1,10] = 2
$data[1,11] = 0.6666666666666666
$data[1,12] = 0.9381496666666668
$data[1,13] = 2.814449
$data[1,14] = 0.5433136653246862
$data[1,15] = 0.5433136653246862
$data[1,16] = 20.29186238442856
$data[1,17] = 182.626761459857
$data[1,18] = 0.002323702767632241
$data[1,19] = 0.002323702767632241
$data[2,0] = "FAPs"
$data[2,1] = "Col3a1"
$data[2,2] = "Mag"
$data[2,3] = "M2"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 4967.017741
$data[2,7] = 14901.053223
$data[2,8] = 0.9821456064948035
$data[2,9] = 0.9821456064948036
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 0.7885686666666668
$data[2,13] = 2.365706
$data[2,14] = 0.4566863346753138
$data[2,15] = 0.4566863346753137
$data[2,16] = 3916.834557330049
$data[2,17] = 35251.51101597044
$data[2,18] = 0.4485324771475749
$data[2,19] = 0.4485324771475749
$data[3,0] = "FAPs"
$data[3,1] = "Col3a1"
$data[3,2] = "Mag"
$data[3,3] = "sCs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 4967.017741
$data[3,7] = 14901.053223
$data[3,8] = 0.9821456064948035
$data[3,9] = 0.9821456064948036
$data[3,10] = 2
$data[3,11] = 0.6666666666666666
$data[3,12] = 0.9381496666666668
$data[3,13] = 2.814449
$data[3,14] = 0.5433136653246862
$data[3,15] = 0.5433136653246862
$data[3,16] = 4659.806038046569
$data[3,17] = 41938.25434241912
$data[3,18] = 0.5336131293472286
$data[3,19] = 0.5336131293472287
$data[4,0] = "M2"
$data[4,1] = "Col3a1"
$data[4,2] = "Mag"
$data[4,3] = "M2"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 2.249417666666667
$data[4,7] = 6.748253
$data[4,8] = 0.0004447851394313067
$data[4,9] = 0.0004447851394313068
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 0.7885686666666668
$data[4,13] = 2.365706
$data[4,14] = 0.4566863346753138
$data[4,15] = 0.4566863346753137
$data[4,16] = 1.773820290179778
$data[4,17] = 15.964382611618
$data[4,18] = 0.0002031272950449318
$data[4,19] = 0.0002031272950449318
$data[5,0] = "M2"
$data[5,1] = "Col3a1"
$data[5,2] = "Mag"
$data[5,3] = "sCs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 2.249417666666667
$data[5,7] = 6.748253
$data[5,8] = 0.0004447851394313067
$data[5,9] = 0.0004447851394313068
$data[5,10] = 2
$data[5,11] = 0.6666666666666666
$data[5,12] = 0.9381496666666668
$data[5,13] = 2.814449
$data[5,14] = 0.5433136653246862
$data[5,15] = 0.5433136653246862
$data[5,16] = 2.110290434177445
$data[5,17] = 18.992613907597
$data[5,18] = 0.0002416578443863749
$data[5,19] = 0.0002416578443863749
$data[6,0] = "sCs"
$data[6,1] = "Col3a1"
$data[6,2] = "Mag"
$data[6,3] = "M2"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 66.41617433333333
$data[6,7] = 199.248523
$data[6,8] = 0.01313269998680205
$data[6,9] = 0.01313269998680205
$data[6,10] = 3
$data[6,11] = 1
$data[6,12] = 0.7885686666666668
$data[6,13] = 2.365706
$data[6,14] = 0.4566863346753138
$data[6,15] = 0.4566863346753137
$data[6,16] = 52.37371403913756
$data[6,17] = 471.363426352238
$data[6,18] = 0.005997524621363172
$data[6,19] = 0.005997524621363171
$data[7,0] = "sCs"
$data[7,1] = "Col3a1"
$data[7,2] = "Mag"
$data[7,3] = "sCs"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 66.41617433333333
$data[7,7] = 199.248523
$data[7,8] = 0.01313269998680205
$data[7,9] = 0.01313269998680205
$data[7,10] = 2
$data[7,11] = 0.6666666666666666
$data[7,12] = 0.9381496666666668
$data[7,13] = 2.814449
$data[7,14] = 0.5433136653246862
$data[7,15] = 0.5433136653246862
$data[7,16] = 62.30831181209189
$data[7,17] = 560.774806308827
$data[7,18] = 0.007135175365438882
$data[7,19] = 0.007135175365438882

$ws.Range("A2:T9").Value = $data
